$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header cells for the two new columns
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Match the header formatting used by the other header cells (bold, border, centered)
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)

# Data for column I and J (rows 2-36)
$data = @(
    @(9,9),
    @(8,8),
    @(7,7),
    @(5,5),
    @(7,8),
    @(8,8),
    @(8,8),
    @(4,4),
    @(1,1),
    @(5,6),
    @(8,8),
    @(6,6),
    @(6,6),
    @(4,5),
    @(5,6),
    @(7,7),
    @(7,7),
    @(5,6),
    @(7,7),
    @(7,7),
    @(9,9),
    @(7,7),
    @(5,5),
    @(7,8),
    @(7,8),
    @(7,7),
    @(4,5),
    @(5,5),
    @(8,9),
    @(4,4),
    @(6,6),
    @(5,5),
    @(5,5),
    @(5,5),
    @(2,2)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 9).Value = $data[$i][0]
    $ws.Cells.Item($row, 10).Value = $data[$i][1]
}
